$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: a single formula typed directly into the first cell (stays a standalone <f>, not shared)
$ws.Range("F1").Formula = "=A1-2*10+2*9"

# F2:F22: the formula filled down the rest of column F -> one shared-formula group
$ws.Range("F2:F22").Formula = "=A2-2*10+2*9"

# G1:H16: formula filled across two columns and down 16 rows -> one shared-formula group
$ws.Range("G1:H16").Formula = "=B1-2*10+2*9"

# G17:G22: continuation of column G filled separately -> its own shared-formula group
$ws.Range("G17:G22").Formula = "=B17-2*10+2*9"

# H17:H22: continuation of column H filled separately -> its own shared-formula group
$ws.Range("H17:H22").Formula = "=C17-2*10+2*9"

$null = $ws.Range("F1:H22").Select()
